$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helper: find the 1-based index of the paragraph whose text starts
# with $searchText.
# ------------------------------------------------------------------
function Get-ParaIndexByText($searchText) {
    $i = 1
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like ($searchText + "*")) {
            return $i
        }
        $i = $i + 1
    }
    return -1
}

# ------------------------------------------------------------------
# Helper: insert a new bold paragraph containing $etapeText right
# before the paragraph whose text starts with $searchText.
# ------------------------------------------------------------------
function Insert-EtapeBefore($searchText, $etapeText) {
    $idx = Get-ParaIndexByText $searchText
    if ($idx -lt 0) {
        return
    }
    $target = $d.Paragraphs($idx)
    $target.Range.InsertParagraphBefore() | Out-Null

    # Re-fetch the freshly created (empty) paragraph at the same index.
    $newPara = $d.Paragraphs($idx)
    $newPara.Range.Text = $etapeText
    $newPara.Range.Bold = 1
}

# 1) "Etape 1 : " before the "git add nomdufichier" command line
Insert-EtapeBefore "git add nomdufichier" "Etape 1 : "

# 2) "Etape 2 :" before the "git commit -m ..." command line
Insert-EtapeBefore "git commit -m" "Etape 2 :"

# 3) "Etape 3 : " before the "git push" command line
Insert-EtapeBefore "git push" "Etape 3 : "

# ------------------------------------------------------------------
# Cosmetic clean-up: merge the leading runs of the explanatory
# paragraph under "Etape 1" back into a single run (matches the
# canonical OOXML produced by the original author's edit).
# ------------------------------------------------------------------
$oldPrefix = "(fait la liste des fichier à envoyer, pour ajouter un fichier et envoyer plusieurs fichiers à la fois, refaire une ligne git "
$rng = $d.Content
$rng.Find.Execute($oldPrefix, $false, $false, $false, $false, $false, $true, 1, $false, $oldPrefix, 2) | Out-Null
